$d = $word.ActiveDocument

# --- Table 1 (first table): restyle to TableGrid and add a small left indent ---
$t1 = $d.Tables.Item(1)
$t1.Style = "TableGrid"
$t1.Rows.LeftIndent = 0.25   # 5 dxa == 0.25 pt -> <w:tblInd w:w="5" w:type="dxa"/>

# --- Table 3 (last table): restyle to TableGrid and tweak the table indent ---
$t3 = $d.Tables.Item(3)
$t3.Style = "TableGrid"
$t3.Rows.LeftIndent = -7.1   # -142 dxa == -7.1 pt -> <w:tblInd w:w="-142" w:type="dxa"/>

# Drop the explicit "center" vertical alignment on the pin_code / qr_code cells
# in the last row of table 3, restoring them to the (top) default.
$lastRow = $t3.Rows.Item($t3.Rows.Count)
$lastRow.Cells.Item(2).VerticalAlignment = 0
$lastRow.Cells.Item(3).VerticalAlignment = 0

Write-Host "Template table styling updated."
